$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = 0.303030303030303
    "C2" = 0.330508474576271
    "D2" = 0.275641025641026
    "E2" = 0.256637168141593
    "F2" = 0.274223034734918

    "B3" = 0.502164502164502
    "C3" = 0.572033898305085
    "D3" = 0.57051282051282
    "E3" = 0.539823008849557
    "F3" = 0.404936014625229

    "B4" = 0.155844155844156
    "C4" = 0.152542372881356
    "D4" = 0.198717948717949
    "E4" = 0.123893805309735
    "F4" = 0.158135283363803

    "B5" = 0.303030303030303
    "C5" = 0.322033898305085
    "D5" = 0.198717948717949
    "E5" = 0.348082595870207
    "F5" = 0.448811700182815

    "B6" = 0.515151515151515
    "C6" = 0.601694915254237
    "D6" = 0.673076923076923
    "E6" = 0.690265486725664
    "F6" = 0.648994515539305
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
